# Trade #9 closed at 2026-02-17 07:53:11 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#   - Summary sheet: roll the new trade into the aggregate stats
#   - Strategy Status sheet: roll the new trade into the MarketMaking row
#   - All Trades / MarketMaking sheets: append the new closed trade as row 10

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("B3").Value = 1200      # Current Capital
$summary.Range("B4").Value = 0         # Total P&L $   (-> -0)
$summary.Range("B5").Value = 0         # Total P&L %   (-> -0)
$summary.Range("B6").Value = 9         # Total Trades
$summary.Range("B8").Value = 5         # Losing Trades
$summary.Range("B9").Value = 44.44     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")

$status.Range("C4").Value = 100        # Capital
$status.Range("D4").Value = 9          # Trades
$status.Range("E4").Value = 0          # P&L $  (-> -0)
$status.Range("F4").Value = 0          # P&L %  (-> -0)
$status.Range("G4").Value = 44.44      # Win Rate %

# ---------------------------------------------------------------------------
# Helper to append the new trade row (row 10) to a trades sheet
# ---------------------------------------------------------------------------
function Add-TradeRow($sheet) {
    # B10/C10 hold date-/time-looking text ("2026-02-17", "07:53:04"). Force
    # Text format before assignment so they aren't auto-converted to date
    # serials, then clear the formatting again so no stray style is left
    # behind on the cell (other text columns, e.g. "CLOSED"/"UP", aren't
    # date-ambiguous and don't need this).
    $sheet.Range("B10:C10").NumberFormat = "@"

    $sheet.Range("A10").Value = 9
    $sheet.Range("B10").Value = "2026-02-17"
    $sheet.Range("C10").Value = "07:53:04"
    $sheet.Range("D10").Value = "MarketMaking"
    $sheet.Range("E10").Value = "UP"
    $sheet.Range("F10").Value = 0.78
    $sheet.Range("G10").Value = 0.74
    $sheet.Range("H10").Value = "CLOSED"
    $sheet.Range("I10").Value = -5.1282
    $sheet.Range("J10").Value = -0.04
    $sheet.Range("K10").Value = 100
    $sheet.Range("L10").Value = 0
    $sheet.Range("M10").Value = 0
    $sheet.Range("N10").Value = 0.6
    $sheet.Range("O10").Value = "Normal spread capture: 19600 bps"
    $sheet.Range("P10").Value = "early_exit"
    $sheet.Range("Q10").Value = 0.14

    $sheet.Range("B10:C10").ClearFormats()
}

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
